$wb = $excel.ActiveWorkbook

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 170
$ws.Range("I38").Value = 170
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 510
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -138
$ws.Range("N38").Value = $null

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2941.3044
$ws.Range("I58").Value = 265
$ws.Range("K58").Value = 795
$ws.Range("M58").Value = -645

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2606.4075
$ws.Range("I116").Value = 1892.8667
$ws.Range("K116").Value = 1892.8667
$ws.Range("M116").Value = 1549.1333

# ALC row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 36270.77
$ws.Range("J128").Value = 36270.77
$ws.Range("L128").Value = 36270.77
$ws.Range("N128").Value = -46230.77

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2873.889
$ws.Range("I138").Value = 3479.3572
$ws.Range("J138").Value = 2600.4517
$ws.Range("K138").Value = 10438.0716
$ws.Range("L138").Value = 7801.355100000001
$ws.Range("M138").Value = -5298.071599999999
$ws.Range("N138").Value = -18081.3551

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9870.9375
$ws.Range("I32").Value = 11538.583
$ws.Range("J32").Value = 4868
$ws.Range("K32").Value = 11538.583
$ws.Range("L32").Value = 4868
$ws.Range("M32").Value = -11251.583
$ws.Range("N32").Value = -5442

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3850
$ws.Range("I74").Value = 746.875
$ws.Range("J74").Value = 6008.696
$ws.Range("K74").Value = 746.875
$ws.Range("L74").Value = 6008.696
$ws.Range("M74").Value = 127.125
$ws.Range("N74").Value = -7756.696

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3850
$ws.Range("I77").Value = 746.875
$ws.Range("J77").Value = 6008.696
$ws.Range("K77").Value = 3734.375
$ws.Range("L77").Value = 30043.48
$ws.Range("M77").Value = 633.625
$ws.Range("N77").Value = -38779.48

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1898.8
$ws.Range("I122").Value = 1785.4445
$ws.Range("J122").Value = 2068.8333
$ws.Range("K122").Value = 5356.333500000001
$ws.Range("L122").Value = 6206.499899999999
$ws.Range("M122").Value = -2906.333500000001
$ws.Range("N122").Value = -11106.4999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2158635.2
$ws.Range("I132").Value = 2655859
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 7967577
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -7965047
$ws.Range("N132").Value = -17057.9999

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2776
$ws.Range("I105").Value = 3095
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 3095
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = -1348
$ws.Range("N105").Value = -4994

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1450.8108
$ws.Range("I31").Value = 1074.6666
$ws.Range("J31").Value = 3062.8572
$ws.Range("K31").Value = 1074.6666
$ws.Range("L31").Value = 3062.8572
$ws.Range("M31").Value = -779.6666
$ws.Range("N31").Value = -3652.8572

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1450.8108
$ws.Range("I34").Value = 1074.6666
$ws.Range("J34").Value = 3062.8572
$ws.Range("K34").Value = 1074.6666
$ws.Range("L34").Value = 3062.8572
$ws.Range("M34").Value = -872.6666
$ws.Range("N34").Value = -3466.8572

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 835.5
$ws.Range("I105").Value = 814.9375
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 814.9375
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 932.0625
$ws.Range("N105").Value = -4494

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 29413414
$ws.Range("I34").Value = 198.33333
$ws.Range("J34").Value = 35716244
$ws.Range("K34").Value = 594.99999
$ws.Range("L34").Value = 107148732
$ws.Range("M34").Value = -510.99999
$ws.Range("N34").Value = -107148900

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 563.5
$ws.Range("I113").Value = 536.913
$ws.Range("J113").Value = 592.619
$ws.Range("K113").Value = 1610.739
$ws.Range("L113").Value = 1777.857
$ws.Range("M113").Value = 559.261
$ws.Range("N113").Value = -6117.857

# CUL row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1860
$ws.Range("I118").Value = 466.66666
$ws.Range("K118").Value = 1399.99998
$ws.Range("M118").Value = -156.9999800000001

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1854494.2
$ws.Range("J131").Value = 2175867.8
$ws.Range("L131").Value = 6527603.399999999
$ws.Range("N131").Value = -6537683.399999999

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 3840.9033
$ws.Range("I134").Value = 2158.4
$ws.Range("J134").Value = 6900
$ws.Range("K134").Value = 6475.200000000001
$ws.Range("L134").Value = 20700
$ws.Range("M134").Value = -1405.200000000001
$ws.Range("N134").Value = -30840

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 122560.6
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 135845.11
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 135845.11
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -137841.11

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 122560.6
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 135845.11
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 679225.5499999999
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -689209.5499999999

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1962.5
$ws.Range("I102").Value = 1450
$ws.Range("J102").Value = 3500
$ws.Range("K102").Value = 1450
$ws.Range("L102").Value = 3500
$ws.Range("M102").Value = 172
$ws.Range("N102").Value = -6744

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 16992.666
$ws.Range("J123").Value = 16992.666
$ws.Range("L123").Value = 16992.666
$ws.Range("N123").Value = -21892.666

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1529.9
$ws.Range("I126").Value = 1071.2858
$ws.Range("J126").Value = 2600
$ws.Range("K126").Value = 3213.8574
$ws.Range("L126").Value = 7800
$ws.Range("M126").Value = -743.8574000000003
$ws.Range("N126").Value = -12740

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2852.5
$ws.Range("I132").Value = 2268.6086
$ws.Range("J132").Value = 4771
$ws.Range("K132").Value = 6805.825800000001
$ws.Range("L132").Value = 14313
$ws.Range("M132").Value = -4275.825800000001
$ws.Range("N132").Value = -19373

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1646.7059
$ws.Range("I7").Value = 1560.6923
$ws.Range("J7").Value = 1926.25
$ws.Range("K7").Value = 1560.6923
$ws.Range("L7").Value = 1926.25
$ws.Range("M7").Value = -1448.6923
$ws.Range("N7").Value = -2150.25

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 378.75
$ws.Range("I22").Value = 330
$ws.Range("J22").Value = 408
$ws.Range("K22").Value = 330
$ws.Range("L22").Value = 408
$ws.Range("M22").Value = -35
$ws.Range("N22").Value = -998

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 378.75
$ws.Range("I27").Value = 330
$ws.Range("J27").Value = 408
$ws.Range("K27").Value = 330
$ws.Range("L27").Value = 408
$ws.Range("M27").Value = -223
$ws.Range("N27").Value = -622

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1377.25
$ws.Range("I40").Value = 1279.3334
$ws.Range("J40").Value = 1671
$ws.Range("K40").Value = 1279.3334
$ws.Range("L40").Value = 1671
$ws.Range("M40").Value = -1143.3334
$ws.Range("N40").Value = -1943

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2779.4
$ws.Range("I61").Value = 1623.5
$ws.Range("J61").Value = 3550
$ws.Range("K61").Value = 1623.5
$ws.Range("L61").Value = 3550
$ws.Range("M61").Value = -1421.5
$ws.Range("N61").Value = -3954

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1325.2
$ws.Range("I82").Value = 1275
$ws.Range("J82").Value = 1358.6666
$ws.Range("K82").Value = 1275
$ws.Range("L82").Value = 1358.6666
$ws.Range("M82").Value = -914
$ws.Range("N82").Value = -2080.6666

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1325.2
$ws.Range("I85").Value = 1275
$ws.Range("J85").Value = 1358.6666
$ws.Range("K85").Value = 1275
$ws.Range("L85").Value = 1358.6666
$ws.Range("M85").Value = -27
$ws.Range("N85").Value = -3854.6666

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2779.4
$ws.Range("I113").Value = 1623.5
$ws.Range("J113").Value = 3550
$ws.Range("K113").Value = 1623.5
$ws.Range("L113").Value = 3550
$ws.Range("M113").Value = 546.5
$ws.Range("N113").Value = -7890

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 9654.786
$ws.Range("I122").Value = 11307.909
$ws.Range("K122").Value = 33923.727
$ws.Range("M122").Value = -31473.727

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1646.7059
$ws.Range("I126").Value = 1560.6923
$ws.Range("J126").Value = 1926.25
$ws.Range("K126").Value = 4682.0769
$ws.Range("L126").Value = 5778.75
$ws.Range("M126").Value = -2212.0769
$ws.Range("N126").Value = -10718.75

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1242.2069
$ws.Range("I132").Value = 696.05
$ws.Range("J132").Value = 2455.889
$ws.Range("K132").Value = 2088.15
$ws.Range("L132").Value = 7367.667
$ws.Range("M132").Value = 441.8500000000004
$ws.Range("N132").Value = -12427.667

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null
